$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: source "rrfcrfr" -> "sal", amount 242 -> 12, date "4/24/42242" -> "3/23/2025"
$ws.Range("A2").Value = "sal"
$ws.Range("B2").Value = 12
# Leading apostrophe forces Excel to keep the date-looking string as literal
# text instead of auto-converting it into a date serial number.
$ws.Range("C2").Value = "'3/23/2025"

# Row 3: source "fevfe" -> "ww", amount 24424 -> 22, date "2/24/24" -> "3/23/2025"
$ws.Range("A3").Value = "ww"
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = "'3/23/2025"
